$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching style of other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Fill H2:H6 with value 1
$ws.Range("H2:H6").Value = 1
